$d = $word.ActiveDocument

$replacements = @(
    @("2025-12-06 Saturday", "2025-12-07 Sunday"),
    @("29×63=1827", "85×13=1105"),
    @("90×18=1620", "95×18=1710"),
    @("96×13=1248", "83×58=4814"),
    @("35×68=2380", "74×57=4218"),
    @("82×51=4182", "91×82=7462"),
    @("75×26=1950", "42×51=2142"),
    @("70×93=6510", "40×92=3680"),
    @("21×20=420", "86×96=8256"),
    @("51×56=2856", "23×60=1380"),
    @("95×16=1520", "12×55=660"),
    @("32×20=640", "83×65=5395"),
    @("27×84=2268", "40×88=3520"),
    @("25×66=1650", "72×50=3600"),
    @("84×42=3528", "87×73=6351"),
    @("52×27=1404", "34×62=2108"),
    @("30×33=990", "22×97=2134"),
    @("24×34=816", "69×87=6003"),
    @("31×88=2728", "81×69=5589"),
    @("51×52=2652", "66×32=2112"),
    @("82×83=6806", "58×29=1682"),
    @("31×68=2108", "67×72=4824"),
    @("16×55=880", "36×11=396"),
    @("99×73=7227", "59×25=1475"),
    @("87×49=4263", "76×19=1444"),
    @("33×41=1353", "96×58=5568")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
